# "Generate Report for Handoff"
# The row for b.md (row 3) moves from "Handed back: in sync with en-US" to
# "Ready for handoff" on every sheet, a new handoff package (b.*.xlf) is
# recorded for both locales, and an error note is attached because the
# handback that exists is stale relative to the newly published source.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f0527839dccb5542e257ab9eacfa10c63f069db0/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/288736d0801e7eb3565ed7f6ac495c67ccca7b8e/e2e/b.md."

# ---------------------------------------------------------------------
# Overview sheet - row 3 is the b.md file
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-07 05:42:55"

# ---------------------------------------------------------------------
# zh-cn sheet - row 3 (Source File Name = b.md)
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-09-07 05:42:44"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.15

# ---------------------------------------------------------------------
# de-de sheet - row 3 (Source File Name = b.md)
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-09-07 05:42:55"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.15
